$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://test17.cliotest.com/backoffice/control/main"
$ws.Range("G2").Value = "https://test17.cliotest.com/cabicentral/control/main"
$ws.Range("J2").Value = "https://test17.cliotest.com/warehouse/control/main"
